$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$s = $ws.Range("H18").Style
"style type: $($s.GetType().FullName)"
"style tostring: $s"
